# Refresh cryptos table cell values (price + 1h volume %) to latest scrape.
# Rows 31/32 and 47/49 additionally swapped rank order (new coin data replaces old).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.770.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.12%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.679.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.53%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D5").Value = "'600.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.76%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'156.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.49%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +6.20%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.130"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.46%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +0.26%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'5.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.40%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -0.03%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'29.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.39%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  -1.18%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.159.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.61%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'65.647.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.13%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.666.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.25%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'12.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.26%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  -0.98%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +1.56%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'352.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.95%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.11%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'69.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.03%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +6.39%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'9.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.48%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'1.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.45%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -2.45%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -5.47%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.44%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  +0.10%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "'Bittensor"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'529.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.33%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "'PancakeSwap"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'2.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.37%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.63%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -3.30%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'5.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.33%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -1.60%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'20.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.72%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.06%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'157.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.90%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'1.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.92%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.04%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'164.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.92%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'4.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.57%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'2.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.18%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.0612"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.27%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'22.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.26%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'Mantle"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.643"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.35%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -2.45%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'BabyDogeCoin"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.0₆0265"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +16.52%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +1.75%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'20.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.23%  "
$ws.Range("E51").Style = "Normal"
